# Question 1's static-evaluation-function description is extended to also
# subtract MIN's piece count, i.e.:
#   "... the number of pieces owned by MAX."
# becomes
#   "... the number of pieces owned by MAX minus the number of players owned by MIN."
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "the number of pieces owned by MAX.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "the number of pieces owned by MAX minus the number of players owned by MIN.",
    2
)
